# Update odds for the Belgrano - Ind. Rivadavia match (row 3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.86
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("R3").Value = 1.53

# Remove the Oriente Petrolero - Santa Cruz match entirely (was row 5).
$ws.Rows.Item(5).Delete()
